$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 3966.794699136231
$ws.Range("D7").Value = 449.6590999096962
